$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "68.878.86"
$ws.Range("E2").Value = "  +1.83%  "
Set-TextValue $ws.Range("D3") "3.737.98"
$ws.Range("E3").Value = "  -1.96%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  +0.32%  "
Set-TextValue $ws.Range("D5") "601.24"
$ws.Range("E5").Value = "  +0.59%  "
Set-TextValue $ws.Range("D6") "167.60"
$ws.Range("E6").Value = "  -5.18%  "
Set-TextValue $ws.Range("D7") "3.741.07"
$ws.Range("E7").Value = "  -2.03%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +2.00%  "
Set-TextValue $ws.Range("D11") "6.34"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("E12").Value = "  -0.96%  "
Set-TextValue $ws.Range("D13") "38.04"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("E14").Value = "  -0.85%  "
Set-TextValue $ws.Range("D15") "4.360.36"
$ws.Range("E15").Value = "  -1.96%  "
Set-TextValue $ws.Range("D16") "3.734.04"
$ws.Range("E16").Value = "  -1.39%  "
Set-TextValue $ws.Range("D17") "68.834.62"
$ws.Range("E17").Value = "  +1.70%  "
Set-TextValue $ws.Range("D18") "7.26"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  +0.50%  "
Set-TextValue $ws.Range("D20") "17.28"
$ws.Range("E20").Value = "  +4.59%  "
Set-TextValue $ws.Range("D21") "497.63"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("E22").Value = "  +11.71%  "
Set-TextValue $ws.Range("D23") "0.724"
$ws.Range("E23").Value = "  -1.00%  "
Set-TextValue $ws.Range("D24") "85.17"
$ws.Range("E24").Value = "  +1.02%  "
Set-TextValue $ws.Range("D25") "2.31"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("E26").Value = "  -8.43%  "
Set-TextValue $ws.Range("D27") "12.37"
$ws.Range("E27").Value = "  +0.76%  "
Set-TextValue $ws.Range("D28") "10.13"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -0.40%  "
Set-TextValue $ws.Range("D31") "2.46"
$ws.Range("E31").Value = "  +0.53%  "
Set-TextValue $ws.Range("D32") "7.96"
$ws.Range("E32").Value = "  +2.76%  "
Set-TextValue $ws.Range("D33") "31.71"
$ws.Range("E33").Value = "  -4.34%  "
Set-TextValue $ws.Range("D34") "3.886.92"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E35").Value = "  -0.78%  "
Set-TextValue $ws.Range("D36") "3.666.46"
$ws.Range("E36").Value = "  -2.16%  "
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  +0.49%  "
Set-TextValue $ws.Range("D38") "1.02"
$ws.Range("E38").Value = "  +1.09%  "
Set-TextValue $ws.Range("D39") "5.82"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  -2.67%  "
Set-TextValue $ws.Range("D41") "0.326"
$ws.Range("E41").Value = "  -0.97%  "
Set-TextValue $ws.Range("D42") "435.79"
$ws.Range("E42").Value = "  -4.84%  "
Set-TextValue $ws.Range("D43") "48.95"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -1.44%  "
Set-TextValue $ws.Range("D45") "2.86"
$ws.Range("E45").Value = "  -1.27%  "
Set-TextValue $ws.Range("D46") "8.41"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("E47").Value = "  +0.00%  "
Set-TextValue $ws.Range("D48") "40.67"
$ws.Range("E48").Value = "  -2.37%  "
Set-TextValue $ws.Range("D49") "141.60"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +0.47%  "
Set-TextValue $ws.Range("D51") "2.743.09"
$ws.Range("E51").Value = "  -3.15%  "
